# Update the "last updated" timestamp string (row 1, col A)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 8 de Septiembre de 2020 a las 14:05"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 6486313
$ws.Range("C4").Value = 738
$ws.Range("E4").Value = 2534098
$ws.Range("G4").Value = 52
$ws.Range("H4").Value = 193586

# Row 15 - Iran
$ws.Range("B15").Value = 391112
$ws.Range("C15").Value = 2302
$ws.Range("D15").Value = 337414
$ws.Range("E15").Value = 31156
$ws.Range("G15").Value = 132
$ws.Range("H15").Value = 22542

# Row 37 - Rumania
$ws.Range("B37").Value = 97033
$ws.Range("C37").Value = 1136
$ws.Range("D37").Value = 40838
$ws.Range("E37").Value = 52228
$ws.Range("G37").Value = 41
$ws.Range("H37").Value = 3967

# Row 58 - Nepal
$ws.Range("B58").Value = 48138
$ws.Range("C58").Value = 902
$ws.Range("D58").Value = 32964
$ws.Range("E58").Value = 14868
$ws.Range("G58").Value = 6
$ws.Range("H58").Value = 306

# Row 73 - Estado de Palestina
$ws.Range("B73").Value = 27363
$ws.Range("C73").Value = 584
$ws.Range("E73").Value = 9903
$ws.Range("G73").Value = 6
$ws.Range("H73").Value = 190

# Row 94 - Consejo Danes para los Refugiados
$ws.Range("B94").Value = 10292
$ws.Range("C94").Value = 59
$ws.Range("D94").Value = 9501
$ws.Range("E94").Value = 531

# Row 147 - Malta
$ws.Range("B147").Value = 2099
$ws.Range("C147").Value = 23
$ws.Range("D147").Value = 1729
$ws.Range("E147").Value = 356

# Row 204 - Nueva Caledonia
$ws.Range("D204").Value = 25
$ws.Range("E204").Value = 1
